$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text "Ready for handoff" -> "In Translation" (whole-cell match) on every sheet
# that reports localization status: Overview (zh-cn/de-de columns) and the per-language detail
# sheets (Status column).
$xlWhole = 1
foreach ($ws in @($overview, $zhcn, $dede)) {
    $ws.Cells.Replace("Ready for handoff", "In Translation", $xlWhole) | Out-Null
}

# Shrink the affected columns so they reflect the width of the new, shorter text
# (was sized for "Ready for handoff", now fits "In Translation").
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
